# Adding commits for LCIA
# - "cyclohexane" and "toluene" are merged into a single "hexane" entry
# - Several mass-flow figures in the stream table (rows 10,14,16-18) are updated
# - A new "Θέρμανση" input-side total formula is added at H28
# - View state (selection) is moved to H28

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename "cyclohexane" occurrences to "hexane" (toluene row reused as hexane row) ---
$ws.Range("H9").Value = "hexane"
$ws.Range("A10").Value = "hexane"
$ws.Range("A17").Value = "hexane"
$ws.Range("H20").Value = "hexane"

# --- Updated mass-flow values ---
$ws.Range("D10").Value = 86.18

$ws.Range("C14").Value = 1768.75
$ws.Range("D14").Value = 124.54

$ws.Range("C16").Value = 60.57
$ws.Range("D16").Value = 175.906

$ws.Range("C17").Value = 1.11
$ws.Range("D17").Value = 85.06

$ws.Range("C18").Value = 1.89
$ws.Range("D18").ClearContents()
$ws.Range("E18").Value = 32.634

# --- New formula for the "Θέρμανση" (heating) input fraction ---
$ws.Range("H28").Formula = "=0.942/F16+B15/F16+B16/F16+B18/F16"

# --- Move selection to reflect where the author was last working ---
$ws.Range("H28").Select() | Out-Null

$excel.Calculate() | Out-Null
